# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal").
#
# Inserts two new observation rows (a "Polar King" durazno lot, categories
# "Especial" and "Primera") right above the existing row 128, pushing the
# rest of the price table (old rows 128:233) down to rows 130:235 — which
# also carries the last two existing rows down to the new 234/235 slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 128 downward by two to make room for the new entries.
$ws.Rows("128:129").Insert()

# Common/fixed columns shared by every data row in this sheet.
$mercado  = "Terminal Hortofrutícola Agro Chillán"
$region   = "Ñuble"
$tipo     = "Fruta"
$prod     = "Frutos de hueso (carozo)"
$categoria = "Durazno"
$unidad   = "$/caja 16 kilos empedrada"
$origen   = "Región de O'Higgins"

# New row 128: Polar King / Especial
$ws.Range("A128").Value = 7
$ws.Range("B128").Value = $mercado
$ws.Range("C128").Value = $region
$ws.Range("D128").Value = 44596
$ws.Range("E128").Value = 16
$ws.Range("F128").Value = $tipo
$ws.Range("G128").Value = 100103
$ws.Range("H128").Value = $prod
$ws.Range("I128").Value = 100103004
$ws.Range("J128").Value = $categoria
$ws.Range("K128").Value = "Polar King"
$ws.Range("L128").Value = "Especial"
$ws.Range("M128").Value = 80
$ws.Range("N128").Value = 11000
$ws.Range("O128").Value = 11000
$ws.Range("P128").Value = 11000
$ws.Range("Q128").Value = $unidad
$ws.Range("R128").Value = $origen
$ws.Range("S128").Value = 688
$ws.Range("T128").Value = 16

# New row 129: Polar King / Primera
$ws.Range("A129").Value = 7
$ws.Range("B129").Value = $mercado
$ws.Range("C129").Value = $region
$ws.Range("D129").Value = 44596
$ws.Range("E129").Value = 16
$ws.Range("F129").Value = $tipo
$ws.Range("G129").Value = 100103
$ws.Range("H129").Value = $prod
$ws.Range("I129").Value = 100103004
$ws.Range("J129").Value = $categoria
$ws.Range("K129").Value = "Polar King"
$ws.Range("L129").Value = "Primera"
$ws.Range("M129").Value = 120
$ws.Range("N129").Value = 9000
$ws.Range("O129").Value = 10000
$ws.Range("P129").Value = 9500
$ws.Range("Q129").Value = $unidad
$ws.Range("R129").Value = $origen
$ws.Range("S129").Value = 594
$ws.Range("T129").Value = 16
